$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-20 Thursday" "2025-02-21 Friday"

Replace-Text "675×8=" "970×4="
Replace-Text "702×5=" "493×7="
Replace-Text "947×7=" "171×2="
Replace-Text "195×2=" "367×3="
Replace-Text "831×3=" "572×7="
Replace-Text "496×6=" "987×4="
Replace-Text "438×5=" "981×3="
Replace-Text "523×8=" "117×9="
Replace-Text "935×8=" "588×5="
Replace-Text "510×8=" "218×5="
Replace-Text "468×2=" "471×7="
Replace-Text "356×6=" "491×8="
Replace-Text "226×7=" "548×3="
Replace-Text "420×5=" "210×7="
Replace-Text "822×6=" "938×5="
Replace-Text "316×9=" "527×3="
Replace-Text "462×7=" "693×4="
Replace-Text "944×8=" "276×8="
Replace-Text "391×3=" "187×2="
Replace-Text "891×4=" "336×9="
Replace-Text "754×3=" "610×6="
Replace-Text "528×2=" "523×4="
Replace-Text "301×5=" "379×3="
Replace-Text "235×3=" "712×9="
Replace-Text "259×4=" "596×8="
